$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the added "Percent Changed" column
$ws.Range("F5").Value = "Percent Changed"
$ws.Columns("F").AutoFit() | Out-Null

# Updated Part B (E) values
$ws.Range("E6").Value = 398250
$ws.Range("E7").Value = 253769
$ws.Range("E8").Value = 49944
$ws.Range("E9").Value = 94537

# New "percent changed" formulas in column F
$ws.Range("F6").Formula = "=(E6-D6)/D6"
$ws.Range("F7").Formula = "=(E7-D7)/D7"
$ws.Range("F8").Formula = "=(E8-D8)/D8"
$ws.Range("F9").Formula = "=(E9-D9)/D9"

# Format the new column as percent with 2 decimals
$ws.Range("F6:F9").NumberFormat = "0.00%"

# Match the recorded selection left after the edit
$ws.Range("O9").Select() | Out-Null
